# Experiment order generation script re-run: regenerates each task-order
# sheet's name and stimulus-file rows, in place (sheet positions / physical
# rIds are unchanged - only names + data are refreshed).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1 (was GNG_TO, A1:B5) -> RS_TO, A1:B3 ("eyes open"/"eyes closed")
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "RS_TO-1651588989776852"
$ws1.Range("B2").Value = "eyes open"
$ws1.Range("B3").Value = "eyes closed"
$ws1.Range("A4:B5").EntireRow.Delete()

# ---------------------------------------------------------------------
# Sheet 2 (was NB_TO, A1:B10) -> GNG_TO, A1:B5 (go_stims/GNG_stims)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "GNG_TO-16515889898671927"
$ws2.Range("B2").Value = "go_stims-1651588989791852.csv"
$ws2.Range("B3").Value = "GNG_stims-16515889898221924.csv"
$ws2.Range("B4").Value = "go_stims-16515889898241925.csv"
$ws2.Range("B5").Value = "GNG_stims-16515889898561924.csv"
$ws2.Range("A6:B10").EntireRow.Delete()

# ---------------------------------------------------------------------
# Sheet 3 (was RS_TO, A1:B3) -> vSAT_TO, A1:B5 (SAT_stims/vSAT_stims)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "vSAT_TO-1651588989963192"
$ws3.Range("B2").Value = "SAT_stims-1651588989884192.csv"
$ws3.Range("B3").Value = "vSAT_stims-16515889899481926.csv"
$ws3.Range("A2").Copy()
$ws3.Range("A4:A5").PasteSpecial(-4122)
$ws3.Range("A4").Value = 2
$ws3.Range("B4").Value = "vSAT_stims-16515889899181929.csv"
$ws3.Range("A5").Value = 3
$ws3.Range("B5").Value = "SAT_stims-1651588989903192.csv"

# ---------------------------------------------------------------------
# Sheet 4 (was TOL_TO, A1:B7) -> TOL_TO, A1:B7 (MM_stims/ZM_stims refresh)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-1651588990059191"
$ws4.Range("B2").Value = "MM_stims-16515889899951918.csv"
$ws4.Range("B3").Value = "ZM_stims-16515889899771914.csv"
$ws4.Range("B4").Value = "MM_stims-16515889900281916.csv"
$ws4.Range("B5").Value = "ZM_stims-16515889899971912.csv"
$ws4.Range("B6").Value = "MM_stims-16515889900571918.csv"
$ws4.Range("B7").Value = "ZM_stims-16515889900321915.csv"

# ---------------------------------------------------------------------
# Sheet 5 (was vSAT_TO, A1:B5) -> NB_TO, A1:B10 (ZB-match/OB/TB)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "NB_TO-16515889929709778"
$ws5.Range("B2").Value = "ZB-match_5-16515889916066494.csv"
$ws5.Range("B3").Value = "OB-16515889917965453.csv"
$ws5.Range("B4").Value = "ZB-match_8-16515889915597708.csv"
$ws5.Range("B5").Value = "OB-16515889920739543.csv"
$ws5.Range("A2").Copy()
$ws5.Range("A6:A10").PasteSpecial(-4122)
$ws5.Range("A6").Value = 4
$ws5.Range("B6").Value = "TB-1651588992843991.csv"
$ws5.Range("A7").Value = 5
$ws5.Range("B7").Value = "OB-16515889921291208.csv"
$ws5.Range("A8").Value = 6
$ws5.Range("B8").Value = "ZB-match_2-1651588991257573.csv"
$ws5.Range("A9").Value = 7
$ws5.Range("B9").Value = "TB-16515889929553173.csv"
$ws5.Range("A10").Value = 8
$ws5.Range("B10").Value = "TB-16515889928127346.csv"
